# CMP73010 (1).docx edit:
#  1. Paragraph ">> >  your stuff after this line >>>" (split across 3 runs by
#     two proofErr marks) becomes a single run ">>>  your stuff after this
#     line >>>" with no proofErr marks.
#  2. Paragraph "Baz changes" (two runs split by a _GoBack bookmark) becomes a
#     short write-up about Git, split into four runs before the (preserved)
#     bookmark, followed by a single space after it.
#
# Both paragraphs are rebuilt with Range.InsertXML so the resulting run
# layout (and absence of any leftover run formatting) matches exactly.

$d = $word.ActiveDocument

function New-PackageXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Paragraph: ">>>  your stuff after this line >>>" -----------------
$quotePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ">>*your*stuff after this line*") {
        $quotePara = $p
        break
    }
}

$quoteBody = '<w:p><w:r><w:t>' +
    '&gt;&gt;&gt;  your stuff after this line &gt;&gt;&gt;' +
    '</w:t></w:r></w:p>'
$quotePara.Range.InsertXML((New-PackageXml $quoteBody))

# --- Paragraph: "Baz changes" -> Git write-up ---------------------------
$gitPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Baz chan*") {
        $gitPara = $p
        break
    }
}

$gitBody = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Git is one the version </w:t></w:r>' +
    '<w:r><w:t>control system that tracks the changes</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> in computer files and</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> also performs coordination of work between multiple user of the same project.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$gitPara.Range.InsertXML((New-PackageXml $gitBody))
